# Refresh cryptocurrency Price (column D) and Volume(1h) (column E) values,
# matching the GitHub Actions "Updated cryptos list" data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.098.23'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '''1.826.00'
$ws.Range("E3").Value = '  +0.38%  '
$ws.Range("D4").Value = '''1.007'
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").Value = '''312.32'
$ws.Range("E5").Value = '  +0.43%  '
$ws.Range("D7").Value = '''0.4684'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '''0.3652'
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("D9").Value = '''0.07388'
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = '''0.8803'
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").Value = '''20.24'
$ws.Range("E11").Value = '  -0.19%  '
$ws.Range("D12").Value = '''1.887.55'
$ws.Range("E12").Value = '  +4.75%  '
$ws.Range("D13").Value = '''0.07336'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").Value = '''93.00'
$ws.Range("E15").Value = '  +1.63%  '
$ws.Range("D16").Value = '''6.526'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").Value = '''0.000008719'
$ws.Range("E18").Value = '  -0.21%  '
$ws.Range("D19").Value = '''1.008'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").Value = '''27.466.73'
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '''5.237'
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '''2.080.37'
$ws.Range("E24").Value = '  +2.27%  '
$ws.Range("D25").Value = '''1.880'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '''151.34'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("D28").Value = '''2.138'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '''5.157'
$ws.Range("E29").Value = '  -1.53%  '
$ws.Range("D30").Value = '''116.32'
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").Value = '''0.7429'
$ws.Range("E32").Value = '  -1.66%  '
$ws.Range("E33").Value = '  +0.39%  '
$ws.Range("D34").Value = '''4.509'
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").Value = '''2.945'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = '''1.008'
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").Value = '''2.530'
$ws.Range("E37").Value = '  +6.20%  '
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").Value = '''0.05282'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").Value = '''0.01933'
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("D41").Value = '''7.324'
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").Value = '''2.933'
$ws.Range("E42").Value = '  -1.21%  '
$ws.Range("D43").Value = '''0.5247'
$ws.Range("E43").Value = '  -0.90%  '
$ws.Range("D44").Value = '''0.1640'
$ws.Range("E44").Value = '  -0.78%  '
$ws.Range("D45").Value = '''8.378'
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").Value = '''0.4885'
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").Value = '''10.39'
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("D48").Value = '''1.008'
$ws.Range("E48").Value = '  +0.48%  '
$ws.Range("D49").Value = '''104.37'
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("D50").Value = '''1.648'
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("E51").Value = '  -0.45%  '
